# Pusapati_LabExam03Grading.xlsx - grading points entered for "Generic" and
# "Customer Class" sections (column E - "Points for grading"), matching the
# deducted/awarded points already present in column D ("Total Points").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Generic section (rows 3-6) ---
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 2

# --- Customer Class section (rows 10-14) ---
$ws.Range("E10").Value = 2
$ws.Range("E11").Value = 2
$ws.Range("E12").Value = 2
$ws.Range("E13").Value = 2
$ws.Range("E14").Value = 2

# Recalculate so the SUM() totals (E7, E15, E38, ...) refresh.
$excel.CalculateFull()

# Move the active selection to E15, matching the saved cursor position.
$ws.Activate()
$ws.Range("E15").Select()
